$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 34, shifting rows 34..130 down to 35..131
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with data (copy of old row 34 with a few
# fields changed: Fecha 44519 -> 44526, Volumen 400 -> 450, Origen Región
# Metropolitana -> Región de O'Higgins)
$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(34, 3).Value = "Metropolitana"
$ws.Cells.Item(34, 4).Value = 44526
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 13
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100101
$ws.Cells.Item(34, 8).Value = "Berries"
$ws.Cells.Item(34, 9).Value = 100101001
$ws.Cells.Item(34, 10).Value = "Arándano (blue)"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 450
$ws.Cells.Item(34, 14).Value = 5000
$ws.Cells.Item(34, 15).Value = 5000
$ws.Cells.Item(34, 16).Value = 5000
$ws.Cells.Item(34, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(34, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(34, 19).Value = 2500
$ws.Cells.Item(34, 20).Value = 2
